# Fruta / hortaliza, semanal
# Insert a new week's worth of "Murcott" price rows (2021-09-09, serial 44448)
# directly above the existing Murcott block at row 230, pushing the existing
# rows (Murcott @44167 and Clemenuless @44399) down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the current row 230 (Murcott/Especial row).
$ws.Rows("230:232").Insert()

# Common (unchanged) column values for this market/product block.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$fecha     = 44448
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100102
$producto  = "Cítricos"
$categoriaId = 100102004
$categoria = "Mandarina"
$variedad  = "Murcott"
$unidad    = "$/bandeja 10 kilos"
$origen    = "Provincia de Limarí"
$kgUnidad  = 10

# Row 230: Especial
$r = 230
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 600
$ws.Cells.Item($r, 14).Value = 4500
$ws.Cells.Item($r, 15).Value = 5000
$ws.Cells.Item($r, 16).Value = 4750
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 475
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 231: Primera
$r = 231
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 600
$ws.Cells.Item($r, 14).Value = 3500
$ws.Cells.Item($r, 15).Value = 4000
$ws.Cells.Item($r, 16).Value = 3750
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 375
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 232: Segunda
$r = 232
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 240
$ws.Cells.Item($r, 14).Value = 2500
$ws.Cells.Item($r, 15).Value = 3000
$ws.Cells.Item($r, 16).Value = 2750
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 275
$ws.Cells.Item($r, 20).Value = $kgUnidad
